$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 688.7273
$ws.Range("I28").Value = 638.7143
$ws.Range("K28").Value = 638.7143
$ws.Range("M28").Value = -153.7143

$ws.Range("H39").Value = 8
$ws.Range("I39").Value = 8
$ws.Range("K39").Value = 24
$ws.Range("M39").Value = 272

$ws.Range("H51").Value = 134999.75
$ws.Range("J51").Value = 151428.14
$ws.Range("L51").Value = 151428.14
$ws.Range("N51").Value = -152396.14

$ws.Range("H53").Value = 124.3125
$ws.Range("I53").Value = 91.583336
$ws.Range("J53").Value = 222.5
$ws.Range("K53").Value = 91.583336
$ws.Range("L53").Value = 222.5
$ws.Range("M53").Value = 545.416664
$ws.Range("N53").Value = -1496.5

$ws.Range("H94").Value = 3364.8462
$ws.Range("I94").Value = 3520.25
$ws.Range("K94").Value = 3520.25
$ws.Range("M94").Value = -3069.25

$ws.Range("H138").Value = 2023.2
$ws.Range("J138").Value = 2082.7632
$ws.Range("L138").Value = 6248.2896
$ws.Range("N138").Value = -16528.2896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 11449444
$ws.Range("J8").Value = 9999.5
$ws.Range("L8").Value = 9999.5
$ws.Range("N8").Value = -10287.5

$ws.Range("H32").Value = 2780.9656
$ws.Range("I32").Value = 826
$ws.Range("K32").Value = 826
$ws.Range("M32").Value = -539

$ws.Range("H39").Value = 8000
$ws.Range("I39").Value = 8000
$ws.Range("K39").Value = 8000
$ws.Range("M39").Value = -7480

$ws.Range("H45").Value = 2849.3333
$ws.Range("I45").Value = 1929.2
$ws.Range("J45").Value = 3999.5
$ws.Range("K45").Value = 1929.2
$ws.Range("L45").Value = 3999.5
$ws.Range("M45").Value = -1552.2
$ws.Range("N45").Value = -4753.5

$ws.Range("H61").Value = 9334
$ws.Range("J61").Value = 4999.5
$ws.Range("L61").Value = 4999.5
$ws.Range("N61").Value = -5423.5

$ws.Range("H136").Value = 9334
$ws.Range("J136").Value = 4999.5
$ws.Range("L136").Value = 14998.5
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 18333
$ws.Range("I96").Value = 18333
$ws.Range("K96").Value = 18333
$ws.Range("M96").Value = -15587

$ws.Range("H107").Value = 1103.4445
$ws.Range("I107").Value = 1159.7142
$ws.Range("J107").Value = 906.5
$ws.Range("K107").Value = 1159.7142
$ws.Range("L107").Value = 906.5
$ws.Range("M107").Value = 760.2858000000001
$ws.Range("N107").Value = -4746.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 24950
$ws.Range("I23").Value = 24950
$ws.Range("K23").Value = 24950
$ws.Range("M23").Value = -24710

$ws.Range("H27").Value = 24950
$ws.Range("I27").Value = 24950
$ws.Range("K27").Value = 24950
$ws.Range("M27").Value = -24758

$ws.Range("H58").Value = 2835.5
$ws.Range("I58").Value = 2308.9167
$ws.Range("K58").Value = 2308.9167
$ws.Range("M58").Value = -2105.9167

$ws.Range("H92").Value = 19000
$ws.Range("J92").Value = 19000
$ws.Range("L92").Value = 19000
$ws.Range("N92").Value = -23992

$ws.Range("H132").Value = 9543.777
$ws.Range("I132").Value = 9316.5
$ws.Range("J132").Value = 9998.333000000001
$ws.Range("K132").Value = 27949.5
$ws.Range("L132").Value = 29994.999
$ws.Range("M132").Value = -25419.5
$ws.Range("N132").Value = -35054.999

$ws.Range("H136").Value = 2835.5
$ws.Range("I136").Value = 2308.9167
$ws.Range("K136").Value = 6926.750100000001
$ws.Range("M136").Value = -4376.750100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 141.95833
$ws.Range("I5").Value = 141.95833
$ws.Range("K5").Value = 425.87499
$ws.Range("M5").Value = -313.87499

$ws.Range("H98").Value = 2070.2856
$ws.Range("I98").Value = 2038.6
$ws.Range("K98").Value = 6115.799999999999
$ws.Range("M98").Value = -4617.799999999999

$ws.Range("H106").Value = 12500
$ws.Range("J106").Value = 12500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -39392

$ws.Range("H132").Value = 2098.3333
$ws.Range("I132").Value = 2098.3333
$ws.Range("K132").Value = 18884.9997
$ws.Range("M132").Value = -16354.9997

$ws.Range("H135").Value = 141.95833
$ws.Range("I135").Value = 141.95833
$ws.Range("K135").Value = 1277.62497
$ws.Range("M135").Value = 1257.37503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 26529.5

$ws.Range("H47").Value = 19971.666
$ws.Range("J47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("N47").Value = -11136

$ws.Range("H50").Value = 26529.5

$ws.Range("H53").Value = 30001
$ws.Range("J53").Value = 30001
$ws.Range("L53").Value = 30001
$ws.Range("N53").Value = -31263

$ws.Range("H109").Value = 42000
$ws.Range("J109").Value = 42000
$ws.Range("L109").Value = 42000
$ws.Range("N109").Value = -44080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2641.9524
$ws.Range("I22").Value = 1766.25
$ws.Range("J22").Value = 3180.8462
$ws.Range("K22").Value = 1766.25
$ws.Range("L22").Value = 3180.8462
$ws.Range("M22").Value = -1471.25
$ws.Range("N22").Value = -3770.8462

$ws.Range("H27").Value = 2641.9524
$ws.Range("I27").Value = 1766.25
$ws.Range("J27").Value = 3180.8462
$ws.Range("K27").Value = 1766.25
$ws.Range("L27").Value = 3180.8462
$ws.Range("M27").Value = -1659.25
$ws.Range("N27").Value = -3394.8462

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H46").Value = 3985.5217
$ws.Range("I46").Value = 2378.2666
$ws.Range("K46").Value = 2378.2666
$ws.Range("M46").Value = -2190.2666

$ws.Range("H47").Value = 24285.428
$ws.Range("I47").Value = 19999.5
$ws.Range("J47").Value = 25999.8
$ws.Range("K47").Value = 19999.5
$ws.Range("L47").Value = 25999.8
$ws.Range("M47").Value = -19509.5
$ws.Range("N47").Value = -26979.8

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H52").Value = 24285.428
$ws.Range("I52").Value = 19999.5
$ws.Range("J52").Value = 25999.8
$ws.Range("K52").Value = 19999.5
$ws.Range("L52").Value = 25999.8
$ws.Range("M52").Value = -19766.5
$ws.Range("N52").Value = -26465.8

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H100").Value = 2588
$ws.Range("I100").Value = 2427.4285
$ws.Range("K100").Value = 2427.4285
$ws.Range("M100").Value = -1886.4285

$ws.Range("H122").Value = 5708.7144
$ws.Range("I122").Value = 3618
$ws.Range("J122").Value = 6995.3076
$ws.Range("K122").Value = 10854
$ws.Range("L122").Value = 20985.9228
$ws.Range("M122").Value = -8404
$ws.Range("N122").Value = -25885.9228

$ws.Range("H136").Value = 7169.9
$ws.Range("I136").Value = 7133.222
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 21399.666
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -18849.666
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 32499.75
$ws.Range("J41").Value = 19999.5
$ws.Range("L41").Value = 19999.5
$ws.Range("N41").Value = -20779.5

$ws.Range("H126").Value = 1372
$ws.Range("I126").Value = 1372
$ws.Range("K126").Value = 4116
$ws.Range("M126").Value = -1646
